$wb = $excel.ActiveWorkbook

# 1. Duplicate current first sheet ("Nadal 10") and insert the copy before it.
$first = $wb.Worksheets.Item(1)
$first.Copy($first)
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "Nädal 11"

# 2. Fill in the new week's header + first log row; clear the rest of the carried-over rows.
$newSheet.Range("G4").Value = "6-12aprill"
$newSheet.Range("B7").Value = 43927
$newSheet.Range("C7").Value = 0.79513888888888884
$newSheet.Range("D7").ClearContents()
$newSheet.Range("F7").ClearContents()
$newSheet.Range("G7").Value = "10. kodutöö"

$newSheet.Range("C8:G8").ClearContents()
$newSheet.Range("B9").ClearContents()
$newSheet.Range("C9:D9").ClearContents()
$newSheet.Range("F9:G9").ClearContents()
$newSheet.Range("C10:G10").ClearContents()
$newSheet.Range("B11").ClearContents()
$newSheet.Range("C11:D11").ClearContents()
$newSheet.Range("F11:G11").ClearContents()
$newSheet.Range("B12:G12").ClearContents()
$newSheet.Range("C13:D13").ClearContents()
$newSheet.Range("F13:G13").ClearContents()
$newSheet.Range("C14:D14").ClearContents()
$newSheet.Range("F14:G14").ClearContents()

# 3. Adjust the new sheet's view (scroll/zoom/selection).
$newSheet.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 107
$win.ScrollRow = 2
$win.ScrollColumn = 1
$newSheet.Range("G11").Select()

# 4. The previous current-week sheet ("Nadal 10", now shifted to position 2) had its
#    selection/zoom touched while reviewing it.
$prevSheet = $wb.Worksheets.Item(2)
$prevSheet.Activate()
$prevSheet.Range("G4:J4").Select()

# 5. The sheet before that ("Nadal 9", now position 3) had its zoom changed too.
$thirdSheet = $wb.Worksheets.Item(3)
$thirdSheet.Activate()
$win3 = $excel.ActiveWindow
$win3.Zoom = 94

# 6. Leave the new current-week sheet active/selected, as it would be after editing it last.
$newSheet.Activate()
